$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A3 value from 1 to 2
$ws.Range("A3").Value = 2

# Duplicate row 3's content into new row 6, with A6 = 5
$ws.Range("A3:S3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial() | Out-Null
$ws.Range("A6").Value = 5

# Duplicate row 5's content into new rows 7 and 8
$ws.Range("A5:S5").Copy() | Out-Null
$ws.Range("A7").PasteSpecial() | Out-Null
$ws.Range("A7").Value = 6

$ws.Range("A5:S5").Copy() | Out-Null
$ws.Range("A8").PasteSpecial() | Out-Null
$ws.Range("A8").Value = 7

# Update view: scroll and selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A8").Select()
